$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 397, pushing existing row 397.. down to 399..
$ws.Rows.Item(397).Insert()
$ws.Rows.Item(397).Insert()

# Make sure the date cells in the two new rows use the same date/time number
# format as the surrounding data (column D elsewhere in the table).
$ws.Range("D397:D398").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 397: Artic Mist / Primera
$ws.Cells.Item(397, 1).Value = 11
$ws.Cells.Item(397, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(397, 3).Value = "Bíobío"
$ws.Cells.Item(397, 4).Value = 45015
$ws.Cells.Item(397, 5).Value = 8
$ws.Cells.Item(397, 6).Value = "Fruta"
$ws.Cells.Item(397, 7).Value = 100103
$ws.Cells.Item(397, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(397, 9).Value = 100103006
$ws.Cells.Item(397, 10).Value = "Nectarín"
$ws.Cells.Item(397, 11).Value = "Artic Mist"
$ws.Cells.Item(397, 12).Value = "Primera"
$ws.Cells.Item(397, 13).Value = 50
$ws.Cells.Item(397, 14).Value = 17000
$ws.Cells.Item(397, 15).Value = 17000
$ws.Cells.Item(397, 16).Value = 17000
$ws.Cells.Item(397, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(397, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(397, 19).Value = 944
$ws.Cells.Item(397, 20).Value = 18

# New row 398: Artic Mist / Segunda
$ws.Cells.Item(398, 1).Value = 11
$ws.Cells.Item(398, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(398, 3).Value = "Bíobío"
$ws.Cells.Item(398, 4).Value = 45015
$ws.Cells.Item(398, 5).Value = 8
$ws.Cells.Item(398, 6).Value = "Fruta"
$ws.Cells.Item(398, 7).Value = 100103
$ws.Cells.Item(398, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(398, 9).Value = 100103006
$ws.Cells.Item(398, 10).Value = "Nectarín"
$ws.Cells.Item(398, 11).Value = "Artic Mist"
$ws.Cells.Item(398, 12).Value = "Segunda"
$ws.Cells.Item(398, 13).Value = 50
$ws.Cells.Item(398, 14).Value = 15000
$ws.Cells.Item(398, 15).Value = 15000
$ws.Cells.Item(398, 16).Value = 15000
$ws.Cells.Item(398, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(398, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(398, 19).Value = 833
$ws.Cells.Item(398, 20).Value = 18
